# Rename the metadata sheet "SwateTemplateMetadata" -> "isa_template" and
# make it the active sheet/tab (it was previously "Harvesting_Isolation"
# that was active).

$wb = $excel.ActiveWorkbook

# Second sheet (rId2 / sheetId 3) is the metadata sheet.
$ws = $wb.Worksheets.Item(2)
$ws.Name = "isa_template"

# Activate it so it becomes the selected tab (workbookView activeTab + the
# worksheet's tabSelected flag move from the first sheet to this one).
$ws.Activate()

# Update the remembered selection on the now-active sheet.
$ws.Range("G12").Select() | Out-Null
